$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values in column C and D
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1

# Update the selection to the full table range
$ws.Range("A1:D5").Select()
